$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Row 23 was missing its Start Time / End Time values - fill them in.
$ws.Range("B23").Value = 0.82361111111111107
$ws.Range("C23").Value = 0.82500000000000007

# Two new daily power records -> grow the table by two rows.
$row24 = $lo.ListRows.Add()
$row25 = $lo.ListRows.Add()

$ws.Range("A24").Value = 43351
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("D24").Formula = "=(C24-B24)* 1440"
$ws.Range("E24").Formula = "=IF(C24>B24, (C24-B24)*1440, (B24-C24)*1440)"
$ws.Range("F24").Formula = "=ABS((C24-B24)*1440)"
$ws.Range("E24").NumberFormat = "General"
$ws.Range("F24").NumberFormat = "General"

$ws.Range("A25").Value = 43352
$ws.Range("B25").Value = 0.35069444444444442
$ws.Range("D25").Formula = "=(C25-B25)* 1440"
$ws.Range("E25").Formula = "=IF(C25>B25, (C25-B25)*1440, (B25-C25)*1440)"
$ws.Range("F25").Formula = "=ABS((C25-B25)*1440)"
$ws.Range("E25").NumberFormat = "General"
$ws.Range("F25").NumberFormat = "General"

# Match the author's final selection/view state.
[void]$ws.Range("C25").Select()
